# Reserve_Level_Plotting_Variables_APA_2020.xlsx
# Update the "Mapping" sheet with the new WGS 84 (EPSG 4269) bounding-box
# coordinates and map-label values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

# Header row (labels stay the same, kept for clarity/idempotency)
$ws.Range("A1").Value = "Res_Bounding_Box"
$ws.Range("B1").Value = "SK_Bounding_Box"
$ws.Range("C1").Value = "Station_Map_Label"
$ws.Range("D1").Value = "SK_WQ_Map_Label"
$ws.Range("E1").Value = "SK_MET_Map_Label"
$ws.Range("F1").Value = "SK_NUT_Map_Label"

# Row 2
$ws.Range("A2").Value = -85.2966
$ws.Range("B2").Value = -85.2609
$ws.Range("C2").Value = "R"
$ws.Range("D2").Value = "R"
$ws.Range("E2").Value = "L"
$ws.Range("F2").Value = "R"

# Row 3
$ws.Range("A3").Value = 29.5247
$ws.Range("B3").Value = 29.5559
$ws.Range("C3").Value = "L"
$ws.Range("D3").Value = "L"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "L"

# Row 4
$ws.Range("A4").Value = -84.6204
$ws.Range("B4").Value = -84.6561
$ws.Range("C4").Value = "L"
$ws.Range("D4").Value = "L"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "L"

# Row 5
$ws.Range("A5").Value = 30.2817
$ws.Range("B5").Value = 30.2506
$ws.Range("C5").Value = "R"
$ws.Range("D5").Value = "R"
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "R"

# Row 6
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "R"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
